$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "data display"
$ws.Range("H10").Value = "line graph(s)"
$ws.Range("H11").Value = "photo(s)"
$ws.Range("H12").Value = "photo(s)"
$ws.Range("H13").Value = "photo(s)"
$ws.Range("H18").Value = "photo(s)"
$ws.Range("H19").Value = "photo(s)"
$ws.Range("H20").Value = "photo(s)"
$ws.Range("H21").Value = "photo(s)"
$ws.Range("H22").Value = "line graph(s)"
$ws.Range("H24").Value = "line graph(s)"
$ws.Range("H25").Value = "line graph(s)"
$ws.Range("H27").Value = "line graph(s)"
$ws.Range("H28").Value = "line graph(s)"
$ws.Range("H29").Value = "line graph(s)"
$ws.Range("H30").Value = "line graph(s)"
$ws.Range("H31").Value = "line graph(s)"
$ws.Range("H35").Value = "line graph(s)"
$ws.Range("H37").Value = "pie chart(s)"
$ws.Range("H40").Value = "bar chart(s)"
$ws.Range("H41").Value = "line graph(s)"
$ws.Range("H42").Value = "line graph(s)"
$ws.Range("H46").Value = "scatter plot(s)"
$ws.Range("H55").Value = "photo(s)"
$ws.Range("H61").Value = "drawing(s)"
$ws.Range("H62").Value = "line graph(s)"

$ws.Columns("I").Delete()
